$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B63").Value = "SingleUseId68"
$ws.Range("C63").Value = "Typographies_button"
$ws.Range("D63").Value = "Center"
$ws.Range("E63").Value = "LTR"
$ws.Range("F63").Value = "DISPATCH"

$ws.Range("B64").Value = "SingleUseId69"
$ws.Range("C64").Value = "Typographies_button"
$ws.Range("D64").Value = "Center"
$ws.Range("E64").Value = "LTR"
$ws.Range("F64").Value = "STOP"

$ws.Range("B65").Value = "SingleUseId71"
$ws.Range("C65").Value = "Typographies_button"
$ws.Range("D65").Value = "Center"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "PAY SALE"

$ws.Range("B66").Value = "SingleUseId72"
$ws.Range("C66").Value = "Typography_label"
$ws.Range("D66").Value = "Center"
$ws.Range("E66").Value = "LTR"
$ws.Range("F66").Value = "Liters filled: <value>"

$ws.Range("B67").Value = "SingleUseId73"
$ws.Range("C67").Value = "Typography_label"
$ws.Range("D67").Value = "Left"
$ws.Range("E67").Value = "LTR"
$ws.Range("F67").Value = "0"

$ws.Range("B68").Value = "SingleUseId75"
$ws.Range("C68").Value = "Typography_label"
$ws.Range("D68").Value = "Left"
$ws.Range("E68").Value = "LTR"
$ws.Range("F68").Value = "Pump selected is: <value>"

$ws.Range("B69").Value = "SingleUseId76"
$ws.Range("C69").Value = "Typography_label"
$ws.Range("D69").Value = "Center"
$ws.Range("E69").Value = "LTR"
$ws.Range("F69").Value = "Regular"

$ws.Range("B70").Value = "SingleUseId77"
$ws.Range("C70").Value = "Typography_label"
$ws.Range("D70").Value = "Center"
$ws.Range("E70").Value = "LTR"
$ws.Range("F70").Value = "Premium"

$ws.Range("B71").Value = "SingleUseId78"
$ws.Range("C71").Value = "Typography_label"
$ws.Range("D71").Value = "Center"
$ws.Range("E71").Value = "LTR"
$ws.Range("F71").Value = "Regular diesel"

$ws.Range("B72").Value = "SingleUseId79"
$ws.Range("C72").Value = "Typography_label"
$ws.Range("D72").Value = "Center"
$ws.Range("E72").Value = "LTR"
$ws.Range("F72").Value = "Premium diesel"

$ws.Range("B73").Value = "SingleUseId80"
$ws.Range("C73").Value = "Typography_label"
$ws.Range("D73").Value = "Left"
$ws.Range("E73").Value = "LTR"
$ws.Range("F73").Value = "1"

$ws.Range("B74").Value = "SingleUseId70"
$ws.Range("C74").Value = "Typographies_button"
$ws.Range("D74").Value = "Center"
$ws.Range("E74").Value = "LTR"
$ws.Range("F74").Value = "CANCEL SALE"
